# Scheduled-runner style market-price refresh for the Sheets workbook.
# For each touched leve row, the currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H:N) are overwritten with freshly pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 23142.143
$ws.Range("J54").Value = 28400
$ws.Range("L54").Value = 28400
$ws.Range("N54").Value = -29372

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3439.825
$ws.Range("I64").Value = 3037.9312
$ws.Range("J64").Value = 4499.364
$ws.Range("K64").Value = 3037.9312
$ws.Range("L64").Value = 4499.364
$ws.Range("M64").Value = -2789.9312
$ws.Range("N64").Value = -4995.364

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3439.825
$ws.Range("I67").Value = 3037.9312
$ws.Range("J67").Value = 4499.364
$ws.Range("K67").Value = 3037.9312
$ws.Range("L67").Value = 4499.364
$ws.Range("M67").Value = -2179.9312
$ws.Range("N67").Value = -6215.364

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 404574.12
$ws.Range("I107").Value = 594335.1
$ws.Range("K107").Value = 594335.1
$ws.Range("M107").Value = -592415.1

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3238.7693
$ws.Range("I137").Value = 5001
$ws.Range("J137").Value = 2918.3635
$ws.Range("K137").Value = 15003
$ws.Range("L137").Value = 8755.0905
$ws.Range("M137").Value = -12453
$ws.Range("N137").Value = -13855.0905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5594.81
$ws.Range("I32").Value = 4597.2817
$ws.Range("J32").Value = 14447.875
$ws.Range("K32").Value = 4597.2817
$ws.Range("L32").Value = 14447.875
$ws.Range("M32").Value = -4310.2817
$ws.Range("N32").Value = -15021.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 30258.334
$ws.Range("J44").Value = 30258.334
$ws.Range("L44").Value = 30258.334
$ws.Range("N44").Value = -31234.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 8689.223
$ws.Range("I74").Value = 11672
$ws.Range("K74").Value = 11672
$ws.Range("M74").Value = -10798

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 8689.223
$ws.Range("I77").Value = 11672
$ws.Range("K77").Value = 58360
$ws.Range("M77").Value = -53992

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 45742
$ws.Range("J140").Value = 45742
$ws.Range("L140").Value = 45742
$ws.Range("N140").Value = -56102

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3209.8171
$ws.Range("I31").Value = 1169.9344
$ws.Range("K31").Value = 1169.9344
$ws.Range("M31").Value = -874.9344000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3209.8171
$ws.Range("I34").Value = 1169.9344
$ws.Range("K34").Value = 1169.9344
$ws.Range("M34").Value = -967.9344000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 188.83333
$ws.Range("I2").Value = 546.5
$ws.Range("J2").Value = 10
$ws.Range("K2").Value = 3279
$ws.Range("L2").Value = 60
$ws.Range("M2").Value = -3166
$ws.Range("N2").Value = -286

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 3030.75
$ws.Range("I63").Value = 1990.25
$ws.Range("J63").Value = 3204.1667
$ws.Range("K63").Value = 5970.75
$ws.Range("L63").Value = 9612.500100000001
$ws.Range("M63").Value = -5221.75
$ws.Range("N63").Value = -11110.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 3030.75
$ws.Range("I66").Value = 1990.25
$ws.Range("J66").Value = 3204.1667
$ws.Range("K66").Value = 17912.25
$ws.Range("L66").Value = 28837.5003
$ws.Range("M66").Value = -14168.25
$ws.Range("N66").Value = -36325.5003

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 9895.833000000001
$ws.Range("I87").Value = 4846.25
$ws.Range("K87").Value = 14538.75
$ws.Range("M87").Value = -13290.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 9895.833000000001
$ws.Range("I90").Value = 4846.25
$ws.Range("K90").Value = 43616.25
$ws.Range("M90").Value = -37376.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1004.63635
$ws.Range("J122").Value = 4116.3335
$ws.Range("L122").Value = 37047.0015
$ws.Range("N122").Value = -41947.0015

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 1348.9474
$ws.Range("I130").Value = 821.6667
$ws.Range("J130").Value = 1592.3077
$ws.Range("K130").Value = 2465.0001
$ws.Range("L130").Value = 4776.9231
$ws.Range("M130").Value = 2554.9999
$ws.Range("N130").Value = -14816.9231

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2000
$ws.Range("I80").Value = 2000
$ws.Range("K80").Value = 2000
$ws.Range("M80").Value = -1002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2000
$ws.Range("I83").Value = 2000
$ws.Range("K83").Value = 10000
$ws.Range("M83").Value = -5008

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1161
$ws.Range("I113").Value = 1237.5
$ws.Range("J113").Value = 549
$ws.Range("K113").Value = 1237.5
$ws.Range("L113").Value = 549
$ws.Range("M113").Value = 932.5
$ws.Range("N113").Value = -4889

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 39044716
$ws.Range("J140").Value = 39044716
$ws.Range("L140").Value = 39044716
$ws.Range("N140").Value = -39055076

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 63613.312
$ws.Range("I46").Value = 250475.5
$ws.Range("J46").Value = 1325.9166
$ws.Range("K46").Value = 250475.5
$ws.Range("L46").Value = 1325.9166
$ws.Range("M46").Value = -250287.5
$ws.Range("N46").Value = -1701.9166

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 20899498
$ws.Range("I136").Value = 103812.1
$ws.Range("K136").Value = 311436.3
$ws.Range("M136").Value = -308886.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 90002
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 90002
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 90002
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -90226

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3027.2727
$ws.Range("I62").Value = 2855.5557
$ws.Range("J62").Value = 3800
$ws.Range("K62").Value = 2855.5557
$ws.Range("L62").Value = 3800
$ws.Range("M62").Value = -2231.5557
$ws.Range("N62").Value = -5048

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 3027.2727
$ws.Range("I65").Value = 2855.5557
$ws.Range("J65").Value = 3800
$ws.Range("K65").Value = 14277.7785
$ws.Range("L65").Value = 19000
$ws.Range("M65").Value = -11157.7785
$ws.Range("N65").Value = -25240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1991.9259
$ws.Range("I81").Value = 1466.8889
$ws.Range("J81").Value = 2254.4443
$ws.Range("K81").Value = 2933.7778
$ws.Range("L81").Value = 4508.8886
$ws.Range("M81").Value = -1872.7778
$ws.Range("N81").Value = -6630.8886

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1991.9259
$ws.Range("I84").Value = 1466.8889
$ws.Range("J84").Value = 2254.4443
$ws.Range("K84").Value = 14668.889
$ws.Range("L84").Value = 22544.443
$ws.Range("M84").Value = -9364.888999999999
$ws.Range("N84").Value = -33152.443

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 459767.2
$ws.Range("I136").Value = 910399.8
$ws.Range("K136").Value = 2731199.4
$ws.Range("M136").Value = -2728649.4

Write-Output "Refreshed 169 cells across 31 leve rows"
